$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44389
$ws.Range("M2").Value = 60
$ws.Range("N2").Value = 17000
$ws.Range("O2").Value = 17000
$ws.Range("P2").Value = 17000
$ws.Range("S2").Value = 1700

# Row 3
$ws.Range("D3").Value = 44389
$ws.Range("L3").Value = 'Segunda'
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 15000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 15000
$ws.Range("S3").Value = 1500

# Row 4
$ws.Range("D4").Value = 44382
$ws.Range("L4").Value = 'Primera'
$ws.Range("M4").Value = 58
$ws.Range("N4").Value = 17000
$ws.Range("O4").Value = 17000
$ws.Range("P4").Value = 17000
$ws.Range("S4").Value = 1700

# Row 5
$ws.Range("D5").Value = 44445
$ws.Range("M5").Value = 68
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 15000
$ws.Range("P5").Value = 15000
$ws.Range("Q5").Value = '$/bandeja 10 kilos'
$ws.Range("S5").Value = 1500
$ws.Range("T5").Value = 10

# Row 6
$ws.Range("D6").Value = 44354
$ws.Range("M6").Value = 45
$ws.Range("N6").Value = 15000
$ws.Range("O6").Value = 15000
$ws.Range("P6").Value = 15000
$ws.Range("S6").Value = 1500

# Row 7
$ws.Range("D7").Value = 44413
$ws.Range("L7").Value = 'Primera'
$ws.Range("M7").Value = 60

# Row 8
$ws.Range("L8").Value = 'Segunda'
$ws.Range("M8").Value = 58
$ws.Range("N8").Value = 13000
$ws.Range("O8").Value = 13000
$ws.Range("P8").Value = 13000
$ws.Range("S8").Value = 1300

# Row 9
$ws.Range("D9").Value = 44323
$ws.Range("L9").Value = 'Primera'
$ws.Range("M9").Value = 48
$ws.Range("N9").Value = 24000
$ws.Range("O9").Value = 24000
$ws.Range("P9").Value = 24000
$ws.Range("Q9").Value = '$/caja 15 kilos granel'
$ws.Range("S9").Value = 1600
$ws.Range("T9").Value = 15

# Row 10
$ws.Range("D10").Value = 44431
$ws.Range("M10").Value = 65
$ws.Range("N10").Value = 18000
$ws.Range("O10").Value = 18000
$ws.Range("P10").Value = 18000
$ws.Range("Q10").Value = '$/bandeja 10 kilos'
$ws.Range("S10").Value = 1800
$ws.Range("T10").Value = 10

# Row 11
$ws.Range("D11").Value = 44431
$ws.Range("L11").Value = 'Segunda'
$ws.Range("M11").Value = 60
$ws.Range("N11").Value = 16000
$ws.Range("O11").Value = 16000
$ws.Range("P11").Value = 16000
$ws.Range("S11").Value = 1600

# Row 13
$ws.Range("D13").Value = 44396
$ws.Range("M13").Value = 60
$ws.Range("N13").Value = 17000
$ws.Range("O13").Value = 17000
$ws.Range("P13").Value = 17000
$ws.Range("S13").Value = 1700

# Row 14
$ws.Range("D14").Value = 44396
$ws.Range("L14").Value = 'Segunda'
$ws.Range("M14").Value = 56
$ws.Range("N14").Value = 15000
$ws.Range("O14").Value = 15000
$ws.Range("P14").Value = 15000
$ws.Range("S14").Value = 1500

# Row 15
$ws.Range("D15").Value = 44385
$ws.Range("L15").Value = 'Primera'
$ws.Range("N15").Value = 17000
$ws.Range("O15").Value = 17000
$ws.Range("P15").Value = 17000
$ws.Range("S15").Value = 1700

# Row 16
$ws.Range("D16").Value = 44385
$ws.Range("L16").Value = 'Segunda'
$ws.Range("M16").Value = 50
$ws.Range("N16").Value = 15000
$ws.Range("O16").Value = 15000
$ws.Range("P16").Value = 15000
$ws.Range("S16").Value = 1500

# Row 17
$ws.Range("D17").Value = 44371
$ws.Range("L17").Value = 'Primera'
$ws.Range("M17").Value = 67
$ws.Range("N17").Value = 18000
$ws.Range("O17").Value = 18000
$ws.Range("P17").Value = 18000
$ws.Range("S17").Value = 1800

# Row 18
$ws.Range("D18").Value = 44370

# Row 19
$ws.Range("D19").Value = 44433
$ws.Range("L19").Value = 'Primera'
$ws.Range("M19").Value = 50
$ws.Range("N19").Value = 18000
$ws.Range("O19").Value = 18000
$ws.Range("P19").Value = 18000
$ws.Range("S19").Value = 1800

# Row 20
$ws.Range("D20").Value = 44321
$ws.Range("M20").Value = 42

# Row 21
$ws.Range("D21").Value = 44398

# Row 22
$ws.Range("D22").Value = 44398

# Row 23
$ws.Range("D23").Value = 44420
$ws.Range("M23").Value = 54
$ws.Range("N23").Value = 18000
$ws.Range("O23").Value = 18000
$ws.Range("P23").Value = 18000
$ws.Range("S23").Value = 1800

# Row 24
$ws.Range("D24").Value = 44420
$ws.Range("M24").Value = 50
$ws.Range("N24").Value = 15000
$ws.Range("O24").Value = 15000
$ws.Range("P24").Value = 15000
$ws.Range("S24").Value = 1500

# Row 25
$ws.Range("D25").Value = 44417
$ws.Range("M25").Value = 56
$ws.Range("N25").Value = 16000
$ws.Range("O25").Value = 16000
$ws.Range("P25").Value = 16000
$ws.Range("S25").Value = 1600

# Row 26
$ws.Range("D26").Value = 44417
$ws.Range("M26").Value = 60
$ws.Range("N26").Value = 14000
$ws.Range("O26").Value = 14000
$ws.Range("P26").Value = 14000
$ws.Range("S26").Value = 1400

# Row 27
$ws.Range("D27").Value = 44441
$ws.Range("M27").Value = 80
$ws.Range("N27").Value = 15000
$ws.Range("O27").Value = 15000
$ws.Range("P27").Value = 15000
$ws.Range("S27").Value = 1500

# Row 28
$ws.Range("D28").Value = 44315

# Row 29
$ws.Range("D29").Value = 44391

# Row 30
$ws.Range("D30").Value = 44391
$ws.Range("L30").Value = 'Segunda'
$ws.Range("M30").Value = 45
$ws.Range("N30").Value = 15000
$ws.Range("O30").Value = 15000
$ws.Range("P30").Value = 15000
$ws.Range("S30").Value = 1500

# Row 31
$ws.Range("D31").Value = 44329
$ws.Range("M31").Value = 50
$ws.Range("N31").Value = 16000
$ws.Range("O31").Value = 16000
$ws.Range("P31").Value = 16000
$ws.Range("S31").Value = 1600

# Row 32
$ws.Range("D32").Value = 44435
$ws.Range("L32").Value = 'Primera'
$ws.Range("M32").Value = 115
$ws.Range("N32").Value = 18000
$ws.Range("O32").Value = 18000
$ws.Range("P32").Value = 18000
$ws.Range("S32").Value = 1800

# Row 33
$ws.Range("A33").Value = 3
$ws.Range("B33").Value = 'Femacal de La Calera'
$ws.Range("C33").Value = 'Coquimbo'
$ws.Range("D33").Value = 44435
$ws.Range("E33").Value = 5
$ws.Range("F33").Value = 'Fruta'
$ws.Range("G33").Value = 100108
$ws.Range("H33").Value = 'Tropicales y subtropicales'
$ws.Range("I33").Value = 100108004
$ws.Range("J33").Value = 'Papaya'
$ws.Range("K33").Value = 'Cultivar IV Región'
$ws.Range("L33").Value = 'Segunda'
$ws.Range("M33").Value = 60
$ws.Range("N33").Value = 16000
$ws.Range("O33").Value = 16000
$ws.Range("P33").Value = 16000
$ws.Range("Q33").Value = '$/bandeja 10 kilos'
$ws.Range("R33").Value = 'Provincia del Elquí'
$ws.Range("S33").Value = 1600
$ws.Range("T33").Value = 10
$ws.Range("D33").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 34
$ws.Range("A34").Value = 3
$ws.Range("B34").Value = 'Femacal de La Calera'
$ws.Range("C34").Value = 'Coquimbo'
$ws.Range("D34").Value = 44319
$ws.Range("E34").Value = 5
$ws.Range("F34").Value = 'Fruta'
$ws.Range("G34").Value = 100108
$ws.Range("H34").Value = 'Tropicales y subtropicales'
$ws.Range("I34").Value = 100108004
$ws.Range("J34").Value = 'Papaya'
$ws.Range("K34").Value = 'Cultivar IV Región'
$ws.Range("L34").Value = 'Primera'
$ws.Range("M34").Value = 60
$ws.Range("N34").Value = 24000
$ws.Range("O34").Value = 24000
$ws.Range("P34").Value = 24000
$ws.Range("Q34").Value = '$/caja 15 kilos granel'
$ws.Range("R34").Value = 'Provincia del Elquí'
$ws.Range("S34").Value = 1600
$ws.Range("T34").Value = 15
$ws.Range("D34").NumberFormat = "YYYY-MM-DD HH:MM:SS"
